# LeetCodeProblems.xlsx edit
# - Adds new rows (19-23, 25-27) to the "Tree's" worksheet with links to
#   additional tree problems (row 24 left blank, matching the author's sheet).
# - Inserts a brand-new "Graph" worksheet between "Tree's" and "CycleSort"
#   containing two links to graph/BFS problems.
# - Updates the view state (active sheet/cell/scroll) left behind by the
#   author's last save.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Extend the "Tree's" worksheet with the newly-added problems.
#    (Populated first so the shared-string table grows in the same order
#    the author's workbook used.)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Tree's")

$ws.Range("A19").Value = "https://leetcode.com/problems/construct-binary-tree-from-preorder-and-inorder-traversal/description/"
$ws.Range("A19").Style = "Hyperlink"

$ws.Range("A20").Value = "https://leetcode.com/problems/serialize-and-deserialize-binary-tree/description/"
$ws.Range("A20").Style = "Hyperlink"

$ws.Range("A21").Value = "https://leetcode.com/problems/path-sum/description/"

$ws.Range("A22").Value = "https://leetcode.com/problems/sum-root-to-leaf-numbers/description/"
$ws.Hyperlinks.Add($ws.Range("A22"), "https://leetcode.com/problems/sum-root-to-leaf-numbers/description/") | Out-Null

$ws.Range("A23").Value = "https://leetcode.com/problems/binary-tree-maximum-path-sum/description/"
$ws.Hyperlinks.Add($ws.Range("A23"), "https://leetcode.com/problems/binary-tree-maximum-path-sum/description/") | Out-Null

# Row 24 intentionally stays empty (matches the source workbook).

$ws.Range("A25").Value = "https://www.geeksforgeeks.org/convert-bst-to-a-binary-tree/"
$ws.Hyperlinks.Add($ws.Range("A25"), "https://www.geeksforgeeks.org/convert-bst-to-a-binary-tree/") | Out-Null

$ws.Range("A26").Value = "https://www.geeksforgeeks.org/reverse-level-order-traversal/"
$ws.Hyperlinks.Add($ws.Range("A26"), "https://www.geeksforgeeks.org/reverse-level-order-traversal/") | Out-Null

$ws.Range("A27").Value = "https://www.geeksforgeeks.org/boundary-traversal-of-binary-tree/"
$ws.Hyperlinks.Add($ws.Range("A27"), "https://www.geeksforgeeks.org/boundary-traversal-of-binary-tree/") | Out-Null

# ---------------------------------------------------------------------
# 2. Insert the new "Graph" worksheet right after "Tree's" (before
#    "CycleSort"), with its two problem links.
# ---------------------------------------------------------------------
$graphSheet = $wb.Worksheets.Add($null, $ws)
$graphSheet.Name = "Graph"

$graphSheet.Range("A1").Value = "https://www.geeksforgeeks.org/minimum-time-required-so-that-all-oranges-become-rotten/"
$graphSheet.Range("A1").Style = "Hyperlink"
$graphSheet.Hyperlinks.Add($graphSheet.Range("A1"), "https://www.geeksforgeeks.org/minimum-time-required-so-that-all-oranges-become-rotten/") | Out-Null

$graphSheet.Range("A2").Value = "https://leetcode.com/problems/detonate-the-maximum-bombs/description/"
$graphSheet.Range("A2").Style = "Hyperlink"
$graphSheet.Hyperlinks.Add($graphSheet.Range("A2"), "https://leetcode.com/problems/detonate-the-maximum-bombs/description/") | Out-Null

$graphSheet.Range("A3").Select() | Out-Null

# ---------------------------------------------------------------------
# 3. Leave the workbook in the same "view" state as the author's last
#    save: Tree's is the active tab, scrolled down to row 13, with A27
#    selected; BinarySearch (no longer the active tab) last had E28
#    selected.
# ---------------------------------------------------------------------
$bsSheet = $wb.Worksheets.Item("BinarySearch")
$bsSheet.Range("E28").Select() | Out-Null

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("A27").Select() | Out-Null
